$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated vm_pu results for the 380 kV case (commit: "case with 380 kV done")
# Data: rows 2-25 correspond to res_bus index 0-23; columns B,C,D,E,F,I,J,K,L,M,N hold voltage magnitudes (p.u.)
$data = @{
    "B2" = "1.02"; "C2" = "1.030642544694266"; "D2" = "1.044306156886103"; "E2" = "1.04042930257216"; "F2" = "1.051684690084117"; "I2" = "1.034952914384028"; "J2" = "1.035782660434934"; "K2" = "1.047077578173522"; "L2" = "1.043211676064769"; "M2" = "1.054435508702294"; "N2" = "1.037253590618683"
    "B3" = "1.02"; "C3" = "1.031720982316364"; "D3" = "1.044908063961058"; "E3" = "1.041340039432241"; "F3" = "1.052568465981809"; "I3" = "1.035040230140196"; "J3" = "1.036501980480662"; "K3" = "1.047491242902842"; "L3" = "1.043932570540555"; "M3" = "1.05513179732754"; "N3" = "1.037973932181383"
    "B4" = "1.02"; "C4" = "1.032419455834217"; "D4" = "1.045297293857799"; "E4" = "1.041930120588173"; "F4" = "1.05314075970318"; "I4" = "1.035095265534745"; "J4" = "1.036967532284899"; "K4" = "1.047757971813062"; "L4" = "1.044399180545315"; "M4" = "1.055582136257027"; "N4" = "1.038440145122582"
    "B5" = "1.02"; "C5" = "1.032713250099852"; "D5" = "1.045460865933381"; "E5" = "1.042178374919821"; "F5" = "1.053381453954182"; "I5" = "1.035118051365076"; "J5" = "1.037163274538273"; "K5" = "1.047869878398807"; "L5" = "1.044595376873837"; "M5" = "1.055771408185176"; "N5" = "1.038636165352408"
    "B6" = "1.02"; "C6" = "1.032762588593325"; "D6" = "1.045488326826387"; "E6" = "1.042220068688818"; "F6" = "1.053421873518082"; "I6" = "1.035121856602268"; "J6" = "1.037196141931905"; "K6" = "1.04788865470323"; "L6" = "1.044628321041791"; "M6" = "1.055803184776503"; "N6" = "1.038669079421509"
    "B7" = "1.02"; "C7" = "1.032423380914785"; "D7" = "1.045299479754479"; "E7" = "1.041933437052083"; "F7" = "1.053143975472585"; "I7" = "1.035095571380301"; "J7" = "1.036970147707756"; "K7" = "1.047759468004924"; "L7" = "1.044401801999664"; "M7" = "1.055584665517909"; "N7" = "1.038442764259639"
    "B8" = "1.02"; "C8" = "1.031006872587745"; "D8" = "1.044509623746641"; "E8" = "1.040736929617138"; "F8" = "1.051983276099613"; "I8" = "1.034982725852348"; "J8" = "1.036025736231226"; "K8" = "1.047217571952792"; "L8" = "1.04345527564055"; "M8" = "1.054670864741075"; "N8" = "1.037497011610487"
    "B9" = "1.02"; "C9" = "1.028515806159424"; "D9" = "1.043116001091382"; "E9" = "1.038634497321932"; "F9" = "1.049941341112887"; "I9" = "1.034772694782719"; "J9" = "1.034362379801828"; "K9" = "1.046255538968038"; "L9" = "1.04178850491329"; "M9" = "1.053059101217389"; "N9" = "1.03583129302417"
    "B10" = "1.02"; "C10" = "1.026858464983114"; "D10" = "1.042185812911755"; "E10" = "1.037236943474344"; "F10" = "1.048582393682214"; "I10" = "1.034625192426621"; "J10" = "1.033254055326025"; "K10" = "1.045609454921895"; "L10" = "1.040678128409117"; "M10" = "1.051983631093532"; "N10" = "1.034721394600485"
    "B11" = "1.02"; "C11" = "1.026141617802582"; "D11" = "1.041782789851109"; "E11" = "1.036632763419406"; "F11" = "1.047994526534114"; "I11" = "1.034559555460087"; "J11" = "1.032774282050333"; "K11" = "1.045328587199913"; "L11" = "1.040197522799414"; "M11" = "1.051517726597053"; "N11" = "1.03424093999171"
    "B12" = "1.02"; "C12" = "1.025875467852431"; "D12" = "1.041633053953051"; "E12" = "1.036408490645931"; "F12" = "1.047776252832217"; "I12" = "1.034534910021946"; "J12" = "1.032596094126688"; "K12" = "1.045224095009674"; "L12" = "1.040019034409274"; "M12" = "1.051344637122635"; "N12" = "1.034062499020775"
    "B13" = "1.02"; "C13" = "1.025932552535493"; "D13" = "1.041665174376834"; "E13" = "1.036456591269195"; "F13" = "1.04782306935486"; "I13" = "1.034540208531899"; "J13" = "1.032634315073589"; "K13" = "1.045246516405036"; "L13" = "1.040057319409514"; "M13" = "1.051381766824186"; "N13" = "1.034100774245805"
    "B14" = "1.02"; "C14" = "1.026119615337173"; "D14" = "1.041770413343178"; "E14" = "1.03661422196635"; "F14" = "1.047976482193239"; "I14" = "1.034557523661859"; "J14" = "1.032759552545399"; "K14" = "1.045319953210243"; "L14" = "1.040182768278892"; "M14" = "1.051503419609584"; "N14" = "1.034226189569189"
    "B15" = "1.02"; "C15" = "1.02623488667175"; "D15" = "1.041835249931137"; "E15" = "1.036711362893342"; "F15" = "1.048071016384524"; "I15" = "1.034568156994772"; "J15" = "1.032836718311293"; "K15" = "1.045365178167059"; "L15" = "1.04026006545416"; "M15" = "1.051578369724122"; "N15" = "1.034303464919323"
    "B16" = "1.02"; "C16" = "1.026906056461786"; "D16" = "1.042212555197869"; "E16" = "1.037277061425526"; "F16" = "1.048621420524879"; "I16" = "1.03462951134824"; "J16" = "1.033285899252765"; "K16" = "1.045628071924374"; "L16" = "1.040710028781143"; "M16" = "1.052014547121923"; "N16" = "1.034753283749253"
    "B17" = "1.02"; "C17" = "1.027327275671314"; "D17" = "1.042449164304216"; "E17" = "1.037632169475931"; "F17" = "1.048966826959877"; "I17" = "1.03466752453198"; "J17" = "1.033567695975488"; "K17" = "1.04579268216105"; "L17" = "1.040992331502372"; "M17" = "1.052288091852397"; "N17" = "1.035035480655645"
    "B18" = "1.02"; "C18" = "1.027573042176202"; "D18" = "1.042587150598642"; "E18" = "1.03783939152464"; "F18" = "1.049168351219555"; "I18" = "1.034689526403453"; "J18" = "1.033732076574673"; "K18" = "1.04588858943718"; "L18" = "1.041157012648146"; "M18" = "1.052447624678941"; "N18" = "1.035200094694135"
    "B19" = "1.02"; "C19" = "1.027656855174909"; "D19" = "1.04263419627187"; "E19" = "1.03791006472613"; "F19" = "1.049237074994284"; "I19" = "1.034696999522325"; "J19" = "1.033788128360737"; "K19" = "1.045921273138641"; "L19" = "1.041213167822048"; "M19" = "1.05250201762367"; "N19" = "1.035256226080165"
    "B20" = "1.02"; "C20" = "1.027282074927503"; "D20" = "1.042423780835293"; "E20" = "1.037594060062151"; "F20" = "1.048929762468655"; "I20" = "1.03466346371622"; "J20" = "1.033537460471292"; "K20" = "1.045775032108543"; "L20" = "1.0409620411483"; "M20" = "1.052258745275233"; "N20" = "1.035005202213565"
    "B21" = "1.02"; "C21" = "1.026064526719607"; "D21" = "1.041739424048475"; "E21" = "1.036567799601"; "F21" = "1.047931303542535"; "I21" = "1.034552432096683"; "J21" = "1.0327226726419"; "K21" = "1.045298332457958"; "L21" = "1.040145825881892"; "M21" = "1.051467596765587"; "N21" = "1.034189257291998"
    "B22" = "1.02"; "C22" = "1.025299694373803"; "D22" = "1.041308939013362"; "E22" = "1.035923397531643"; "F22" = "1.04730403291236"; "I22" = "1.034481089507712"; "J22" = "1.032210506251158"; "K22" = "1.044997656038098"; "L22" = "1.039632812462097"; "M22" = "1.050969986844684"; "N22" = "1.033676363566234"
    "B23" = "1.02"; "C23" = "1.025705081221456"; "D23" = "1.041537165933673"; "E23" = "1.036264926531102"; "F23" = "1.047636513105881"; "I23" = "1.034519054638107"; "J23" = "1.032482003540328"; "K23" = "1.045157140552003"; "L23" = "1.039904753909906"; "M23" = "1.051233796308429"; "N23" = "1.033948246412697"
    "B24" = "1.02"; "C24" = "1.027302498961579"; "D24" = "1.042435250605976"; "E24" = "1.037611279776878"; "F24" = "1.04894651014842"; "I24" = "1.034665299151171"; "J24" = "1.033551122552923"; "K24" = "1.045783007738534"; "L24" = "1.040975727997565"; "M24" = "1.052272005795776"; "N24" = "1.035018883696919"
    "B25" = "1.02"; "C25" = "1.029159213100506"; "D25" = "1.043476488318709"; "E25" = "1.039177313613978"; "F25" = "1.050468823191771"; "I25" = "1.034828314248869"; "J25" = "1.034792297343534"; "K25" = "1.046505086410887"; "L25" = "1.042219266651336"; "M25" = "1.053475955058337"; "N25" = "1.036261821098097"
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = [double]$data[$addr]
}
